$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Phone"
$ws.Range("C1").Value = "Email"
$ws.Range("A1:C1").Font.Bold = $true

# Data rows
$ws.Range("A2").Value = "John Doe"
$ws.Range("B2").Value = "123-456-7890"
$ws.Range("C2").Value = "john@example.com"

$ws.Range("A3").Value = "Jane Smith"
$ws.Range("B3").Value = "987-654-3210"
$ws.Range("C3").Value = "jane@example.com"
